$wb = $excel.ActiveWorkbook

# --- Aufgaben: zeiten in format hh:mm:ss angepasst -------------------------
# A handful of duration cells in column C were entered an order of magnitude
# too small (minutes stored where hours:minutes were meant). Re-scale them by
# 60 (i.e. reinterpret "seconds of a day" as "minutes of a day") so they read
# correctly under the existing h:mm:ss number format.
$ws = $wb.Worksheets.Item("Aufgaben")
$rows = @(23, 24, 25, 26, 30, 31, 32, 33)
foreach ($r in $rows) {
    $cell = $ws.Range("C$r")
    $old = $cell.Value()
    $secondsOfDay = [math]::Round($old * 86400)
    $new = ($secondsOfDay * 60) / 86400
    $cell.Value = $new
}

# --- Date-format cleanup: use the built-in date format instead of the -----
# --- custom "m/d/yyyy" one for the (currently empty) date cells -----------
$ws.Range("B5:C5").NumberFormat = "mm-dd-yy"
$wsFragen = $wb.Worksheets.Item("Fragen")
$wsFragen.Range("B5:C5").NumberFormat = "mm-dd-yy"

# --- View state on the Aufgaben sheet --------------------------------------
$ws.Activate()
$win = $excel.ActiveWindow
$win.Zoom = 85
[void]$ws.Range("C34").Select()
